$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.788.31"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "3.126.99"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.47%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "3.125.24"
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.475"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.38%  "
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.108"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.414"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.97%  "
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("D14").Value = "3.660.29"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.05%  "
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("D17").Value = "57.853.73"
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("D18").Value = "3.123.60"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("E19").Value = "  +2.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.09%  "
$ws.Range("E21").Value = "  +3.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "368.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.49%  "
$ws.Range("E26").Value = "  +1.47%  "
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").Value = "0.0₃0867"
$ws.Range("E29").Value = "  -2.53%  "
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.16"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.53"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.61%  "
$ws.Range("E38").Value = "  +5.64%  "
$ws.Range("E39").Value = "  -0.96%  "
$ws.Range("E40").Value = "  +4.58%  "
$ws.Range("E41").Value = "  +2.55%  "
$ws.Range("D42").Value = "2.544.50"
$ws.Range("E42").Value = "  +7.07%  "
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("E44").Value = "  +0.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "37.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0270"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.979"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.743"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.68%  "
